$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 8. References and documentation -> "Organization website (if available)" (row 10, col B):
# the statistics committee's site moved from www.stat.kg to www.stat.gov.kg
$ws.Range("B10").Value = "www.stat.gov.kg"

# B2 ("Goal" value cell) should wrap its text (style changes from non-wrapping to wrapping)
$ws.Range("B2").WrapText = $true
